$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Regular price/volume updates ---
$ws.Range("D2").Value = "'30.249.60"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "'1.867.04"
$ws.Range("E3").Value = "  -1.16%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'235.25"
$ws.Range("E5").Value = "  -1.73%  "
$ws.Range("D7").Value = "'0.4661"
$ws.Range("E7").Value = "  -0.44%  "
$ws.Range("D8").Value = "'0.2835"
$ws.Range("E8").Value = "  -1.32%  "
$ws.Range("D9").Value = "'0.06555"
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("D10").Value = "'21.37"
$ws.Range("E10").Value = "  +6.43%  "
$ws.Range("D11").Value = "'0.07873"
$ws.Range("E11").Value = "  +1.17%  "
$ws.Range("D12").Value = "'97.86"
$ws.Range("E12").Value = "  -0.66%  "
$ws.Range("D13").Value = "'1.866.70"
$ws.Range("E13").Value = "  -1.24%  "
$ws.Range("D14").Value = "'5.118"
$ws.Range("E14").Value = "  -0.60%  "
$ws.Range("D15").Value = "'0.6748"
$ws.Range("E15").Value = "  -1.24%  "
$ws.Range("D16").Value = "'281.36"
$ws.Range("E16").Value = "  -1.41%  "
$ws.Range("D17").Value = "'30.241.82"
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("D19").Value = "'5.532"
$ws.Range("E19").Value = "  +2.39%  "
$ws.Range("D20").Value = "'12.69"
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("D21").Value = "'2.112.51"
$ws.Range("E21").Value = "  -1.12%  "
$ws.Range("D22").Value = "'0.000007286"
$ws.Range("E22").Value = "  -0.78%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "'6.176"
$ws.Range("E24").Value = "  -0.68%  "
$ws.Range("D25").Value = "'9.220"
$ws.Range("E25").Value = "  -2.00%  "
$ws.Range("D26").Value = "'164.98"
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("D27").Value = "'19.23"
$ws.Range("E27").Value = "  -0.71%  "
$ws.Range("D29").Value = "'1.375"
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D30").Value = "'0.09692"
$ws.Range("E30").Value = "  -0.94%  "
$ws.Range("D32").Value = "'1.476"
$ws.Range("E32").Value = "  -0.91%  "
$ws.Range("D33").Value = "'4.107"
$ws.Range("E33").Value = "  -1.85%  "
$ws.Range("D34").Value = "'0.04696"
$ws.Range("E34").Value = "  -0.63%  "
$ws.Range("D35").Value = "'1.120"
$ws.Range("E35").Value = "  +1.79%  "
$ws.Range("D37").Value = "'2.729"
$ws.Range("E37").Value = "  +0.83%  "
$ws.Range("D38").Value = "'0.01859"
$ws.Range("E38").Value = "  -1.44%  "
$ws.Range("D41").Value = "'73.85"
$ws.Range("E41").Value = "  +1.06%  "
$ws.Range("D42").Value = "'1.953"
$ws.Range("E42").Value = "  -1.65%  "
$ws.Range("D43").Value = "'0.8459"
$ws.Range("E43").Value = "  -3.08%  "
$ws.Range("D44").Value = "'0.4179"
$ws.Range("E44").Value = "  -1.39%  "
$ws.Range("D47").Value = "'7.202"
$ws.Range("E47").Value = "  -1.44%  "
$ws.Range("D48").Value = "'9.214"
$ws.Range("E48").Value = "  -1.12%  "
$ws.Range("D49").Value = "'935.30"
$ws.Range("E49").Value = "  -5.52%  "
$ws.Range("D50").Value = "'34.08"
$ws.Range("E50").Value = "  -0.45%  "
$ws.Range("D51").Value = "'0.1126"
$ws.Range("E51").Value = "  -3.55%  "

# --- Volume-only updates (price unchanged) ---
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("E28").Value = "  -3.92%  "
$ws.Range("E31").Value = "  -1.39%  "
$ws.Range("E36").Value = "  -1.34%  "

# --- Row 39/40 swap: FraxShare <-> MXToken reorder ---
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.539"
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'6.264"
$ws.Range("E40").Value = "  -6.34%  "

# --- Row 45/46 swap: PaxDollar <-> Quant reorder ---
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'104.01"
$ws.Range("E45").Value = "  -0.44%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = "  +0.11%  "
